# RF013 - Gerenciar Niveis das Competencias workbook update: v1.1 -> v1.2
#
# TC6 and TC7's bodies are swapped:
#   - TC6 previously held the 7-step "Editar" (edit) flow ending in an
#     edit-validation error message (rows 63-69).
#   - TC7 previously held the "TC7" header block (rows 72-75) followed by
#     the 4-step "Excluir" (delete) flow ending in a delete error message
#     (rows 76-79).
# After the edit:
#   - TC6 keeps its own header (rows 59-62, untouched) but now holds the
#     4-step "Excluir" flow (rows 63-66).
#   - Rows 67-68 are blank (same gap that used to sit at rows 70-71).
#   - The "TC7" header block now sits at rows 69-72.
#   - TC7 now holds the 7-step "Editar" flow (rows 73-79).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stage the three blocks that need to move, far away from the
#        working area, so later writes never clobber not-yet-read data.
$ws.Range("A63:F69").Copy()                 # Block E: 7-step "Editar" flow
$ws.Range("A500").PasteSpecial()

$ws.Range("A76:F79").Copy()                 # Block X: 4-step "Excluir" flow
$ws.Range("A510").PasteSpecial()

$ws.Range("A72:F75").Copy()                 # Block H: "TC7" header block
$ws.Range("A520").PasteSpecial()

# --- 2. Remove the stale merged ranges inside the region we're about to
#        rebuild (Clear() does not itself lift merges).
$ws.Range("B73:D73").UnMerge()
$ws.Range("B74:F74").UnMerge()

# --- 3. Wipe the working region (rows 63-79) completely, content+format.
$ws.Range("A63:F79").Clear()

# --- 4. Write the blocks back in their new order:
#        Block X (Excluir) -> 63-66, gap -> 67-68, Block H (TC7 header)
#        -> 69-72, Block E (Editar) -> 73-79.
$ws.Range("A510:F513").Copy()
$ws.Range("A63").PasteSpecial()

$ws.Range("A520:F523").Copy()
$ws.Range("A69").PasteSpecial()

$ws.Range("A500:F506").Copy()
$ws.Range("A73").PasteSpecial()

# --- 5. Clean up the staging area.
$ws.Range("A500:F523").Clear()

$excel.CutCopyMode = $false
